$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to text format so numeric/percent-looking strings
# are preserved exactly as text, matching the source inline-string cells.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "306.23"
$ws.Range("E2").Value = "-0.50%"
$ws.Range("D3").Value = "38.81"
$ws.Range("E3").Value = "7.20%"
$ws.Range("D4").Value = "5.115"
$ws.Range("E4").Value = "1.16%"
$ws.Range("D5").Value = "0.08084"
$ws.Range("E5").Value = "-0.76%"
$ws.Range("D6").Value = "1.927"
$ws.Range("E6").Value = "-3.08%"
$ws.Range("D7").Value = "4.192"
$ws.Range("E7").Value = "0.62%"
$ws.Range("D8").Value = "7.999"
$ws.Range("E8").Value = "1.62%"
$ws.Range("D9").Value = "0.9311"
$ws.Range("E9").Value = "0.26%"
$ws.Range("D10").Value = "0.1463"
$ws.Range("E10").Value = "0.35%"
$ws.Range("D11").Value = "0.1924"
$ws.Range("E11").Value = "-0.54%"
$ws.Range("D12").Value = "0.09092"
$ws.Range("E12").Value = "-0.12%"
$ws.Range("D13").Value = "0.03510"
$ws.Range("E13").Value = "1.72%"
$ws.Range("D14").Value = "0.09799"
$ws.Range("E14").Value = "-0.94%"
$ws.Range("D15").Value = "0.001395"
$ws.Range("E15").Value = "-1.08%"
$ws.Range("D16").Value = "0.005932"
$ws.Range("E16").Value = "-12.57%"
$ws.Range("D17").Value = "3.779"
$ws.Range("E17").Value = "-1.45%"
$ws.Range("E18").Value = "0.52%"
$ws.Range("E19").Value = "-0.08%"
$ws.Range("E20").Value = "4.97%"
$ws.Range("D21").Value = "4.676"
$ws.Range("E21").Value = "-3.40%"
$ws.Range("E22").Value = "3.15%"
$ws.Range("D23").Value = "0.04382"
$ws.Range("E23").Value = "-0.26%"
$ws.Range("E24").Value = "0.32%"
$ws.Range("D25").Value = "0.004272"
$ws.Range("E25").Value = "1.89%"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("E26").Value = "0.15%"
$ws.Range("D39").Value = "0.02033"
$ws.Range("E39").Value = "-0.57%"
$ws.Range("D40").Value = "0.05065"
$ws.Range("E40").Value = "-1.57%"
$ws.Range("D41").Value = "0.007559"
$ws.Range("E41").Value = "1.35%"
$ws.Range("D42").Value = "0.009745"
$ws.Range("E42").Value = "-3.90%"
$ws.Range("E43").Value = "-2.27%"
$ws.Range("D44").Value = "0.002124"
$ws.Range("E44").Value = "-0.32%"
$ws.Range("D45").Value = "0.009936"
$ws.Range("E45").Value = "0.52%"
$ws.Range("D46").Value = "0.00006188"
$ws.Range("E46").Value = "-2.16%"
$ws.Range("D49").Value = "0.001805"
$ws.Range("D50").Value = "0.00002106"
$ws.Range("D51").Value = "0.0002006"

# Restore the default (Normal) cell style so no stray number-format
# style index is left behind on the edited cells.
$priceVolumeRange.Style = "Normal"
